$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 151, shifting the existing rows 151:224 down to 152:225
$ws.Rows.Item(151).Insert()

# Populate the newly inserted row 151 with the latest weekly data point
$ws.Cells.Item(151, 1).Value = 10
$ws.Cells.Item(151, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(151, 3).Value = "La Araucanía"
$ws.Cells.Item(151, 4).Value = 44806
$ws.Cells.Item(151, 5).Value = 9
$ws.Cells.Item(151, 6).Value = "Fruta"
$ws.Cells.Item(151, 7).Value = 100104
$ws.Cells.Item(151, 8).Value = "Frutos de pepita"
$ws.Cells.Item(151, 9).Value = 100104003
$ws.Cells.Item(151, 10).Value = "Membrillo"
$ws.Cells.Item(151, 11).Value = "Champion"
$ws.Cells.Item(151, 12).Value = "Primera"
$ws.Cells.Item(151, 13).Value = 80
$ws.Cells.Item(151, 14).Value = 10000
$ws.Cells.Item(151, 15).Value = 10000
$ws.Cells.Item(151, 16).Value = 10000
$ws.Cells.Item(151, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(151, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(151, 19).Value = 556
$ws.Cells.Item(151, 20).Value = 18
